$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.198.02"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "'1.861.60"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'236.00"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Value = "'0.4674"
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("D8").Value = "'0.2846"
$ws.Range("E8").Value = "  +1.51%  "
$ws.Range("D9").Value = "'0.06525"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").Value = "'21.87"
$ws.Range("E10").Value = "  +12.25%  "
$ws.Range("D11").Value = "'0.07900"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").Value = "'97.20"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("D13").Value = "'1.864.81"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").Value = "'5.154"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").Value = "'0.6788"
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("D16").Value = "'279.11"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "'30.194.39"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "'13.46"
$ws.Range("E18").Value = "  +7.17%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "'0.000007317"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").Value = "'2.110.05"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "'5.358"
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'6.159"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").Value = "'168.38"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("D26").Value = "'9.231"
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("D28").Value = "'1.930"
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").Value = "'1.379"
$ws.Range("E29").Value = "  +3.22%  "
$ws.Range("D30").Value = "'0.09730"
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("D31").Value = "'4.369"
$ws.Range("E31").Value = "  -1.26%  "
$ws.Range("D32").Value = "'1.478"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("D34").Value = "'0.04716"
$ws.Range("E34").Value = "  +1.90%  "
$ws.Range("D35").Value = "'1.134"
$ws.Range("E35").Value = "  +3.17%  "
$ws.Range("D36").Value = "'0.7068"
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "'0.01864"
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("D39").Value = "'2.608"
$ws.Range("E39").Value = "  +4.05%  "
$ws.Range("D40").Value = "'6.293"
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("D41").Value = "'74.67"
$ws.Range("E41").Value = "  +3.76%  "
$ws.Range("D42").Value = "'1.951"
$ws.Range("E42").Value = "  +2.03%  "
$ws.Range("D43").Value = "'0.8478"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").Value = "'0.4171"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").Value = "'103.35"
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("D47").Value = "'967.80"
$ws.Range("E47").Value = "  -2.50%  "
$ws.Range("D48").Value = "'7.183"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").Value = "'9.283"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("D50").Value = "'34.07"
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("E51").Value = "  +0.20%  "
